# Apply "new TPM" data update to the Fgf5-Fgfr2 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared-string label used for the "Inflammatory-Mac" cluster.
# (Row 5's Target cluster still uses this label, now renamed "Resolving-Mac".)
$ws.Range("D5").Value = "Resolving-Mac"

# Row 4's Target cluster becomes "MuSCs" (previously held the Inflammatory-Mac
# / Resolving-Mac label; the data below got swapped with row 5 in the update).
$ws.Range("D4").Value = "MuSCs"

# --- Row 2 (Target cluster: ECs) ---
$ws.Range("G2").Value = 0.2284785
$ws.Range("H2").Value = 0.456957
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2858606666666667
$ws.Range("N2").Value = 0.857582
$ws.Range("O2").Value = 0.0687156860066334
$ws.Range("P2").Value = 0.06932858672617494
$ws.Range("Q2").Value = 0.065313016329
$ws.Range("R2").Value = 0.391878097974
$ws.Range("S2").Value = 0.0687156860066334
$ws.Range("T2").Value = 0.06932858672617494

# --- Row 3 (Target cluster: FAPs) ---
$ws.Range("G3").Value = 0.2284785
$ws.Range("H3").Value = 0.456957
$ws.Range("O3").Value = 0.9046431256549901
$ws.Range("P3").Value = 0.9127119736118995
$ws.Range("Q3").Value = 0.8598469239195
$ws.Range("R3").Value = 5.159081543517001
$ws.Range("S3").Value = 0.9046431256549901
$ws.Range("T3").Value = 0.9127119736118995

# --- Row 4 (Target cluster: MuSCs) ---
$ws.Range("G4").Value = 0.2284785
$ws.Range("H4").Value = 0.456957
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110331
$ws.Range("N4").Value = 0.220662
$ws.Range("O4").Value = 0.02652155835639462
$ws.Range("P4").Value = 0.01783874265571248
$ws.Range("Q4").Value = 0.0252082613835
$ws.Range("R4").Value = 0.100833045534
$ws.Range("S4").Value = 0.02652155835639462
$ws.Range("T4").Value = 0.01783874265571248

# --- Row 5 (Target cluster: Resolving-Mac) ---
$ws.Range("G5").Value = 0.2284785
$ws.Range("H5").Value = 0.456957
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.0004976666666666667
$ws.Range("N5").Value = 0.001493
$ws.Range("O5").Value = 0.0001196299819817856
$ws.Range("P5").Value = 0.0001206970062130259
$ws.Range("Q5").Value = 0.0001137061335
$ws.Range("R5").Value = 0.000682236801
$ws.Range("S5").Value = 0.0001196299819817856
$ws.Range("T5").Value = 0.0001206970062130259
